$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(156, 1).Value = "2023-12-10 12:01:48"
$ws.Cells.Item(156, 2).Value = 0.0004

$ws.Cells.Item(157, 1).Value = "2023-12-10 12:01:58"
$ws.Cells.Item(157, 2).Value = 0.0004
